$d = $word.ActiveDocument

# --- 1. "Skalowanie na przebiegu..." paragraph: remove the (now-relocated)
#        _GoBack bookmark and merge the two runs it used to separate back
#        into a single contiguous run. ---------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$findRange = $d.Content
[void]$findRange.Find.Execute("Skalowanie na przebiegu rejestracyjnym", $true, $false, $false,
                               $false, $false, $true, 1, $false, "", 0)
$skal = $findRange.Paragraphs(1)
$skalStart = $skal.Range.Start
$skalText = $skal.Range.Text
$skalRange = $d.Range($skalStart, $skalStart + $skalText.Length)
$skalRange.Text = "Skalowanie na przebiegu rejestracyjnym nie działa – działa tylko w trakcie rejestracji"

# --- 2. "(OKNO EDYCJI) Przełączanie kursorów..." paragraph: strike the
#        first part of the line and move the _GoBack bookmark here,
#        splitting the trailing "(bez zmiany osi czasu)," into its own
#        un-struck run. -------------------------------------------------
$findRange2 = $d.Content
[void]$findRange2.Find.Execute("Przełączanie kursorów pomiarowych na inne kanały", $true, $false,
                                $false, $false, $false, $true, 1, $false, "", 0)
$przel = $findRange2.Paragraphs(1)
$przelStart = $przel.Range.Start

$part1 = "(OKNO EDYCJI) "
$part2 = "Przełączanie kursorów pomiarowych na inne kanały"
$part3 = " "

# Re-create the _GoBack bookmark spanning from the start of the paragraph
# through the single space that follows "kanały" (before the bookmark is
# split apart by the strikethrough formatting below).
$bmEnd = $przelStart + $part1.Length + $part2.Length + $part3.Length
$bmRange = $d.Range($przelStart, $bmEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)

$r1 = $d.Range($przelStart, $przelStart + $part1.Length)
$r1.Font.StrikeThrough = 1

$r2 = $d.Range($przelStart + $part1.Length, $przelStart + $part1.Length + $part2.Length)
$r2.Font.StrikeThrough = 1
